# Update the occupancy-schedule CSV paths (user profile changed from
# "walkerl\Documents\code" to "LW_Simulation\Documents") and move the
# active-cell selection, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'C:\Users\LW_Simulation\Documents\RC_BuildingSimulator\rc_simulator\auxiliary\occupancy_single_res.csv'
$ws.Range("B3").Value = 'C:\Users\LW_Simulation\Documents\RC_BuildingSimulator\rc_simulator\auxiliary\occupancy_office.csv'

$ws.Range("H6").Select()
